$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 14
$ws.Range("B5").Value = 15
$ws.Range("C6").Value = 14
$ws.Range("B9").Value = 24
$ws.Range("B10").Value = 12

$ws.Range("D11").Font.Size = 11
$ws.Rows("11:11").RowHeight = 15.75

$ws.Range("G14").Select()
